$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A94").Value = 0.853
$ws.Range("B94").Value = 0
$ws.Range("C94").Value = 1.654
